$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.143.55"
$ws.Range("E2").Value = "  +2.17%  "

$ws.Range("D3").Value = "1.916.43"
$ws.Range("E3").Value = "  +2.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.84%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.34"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  -0.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4842"
$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3827"
$ws.Range("E8").Value = "  +1.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07373"
$ws.Range("E9").Value = "  -0.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9379"
$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.84"
$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07824"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").Value = "1.912.00"
$ws.Range("E13").Value = "  +1.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.501"
$ws.Range("E14").Value = "  +0.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.643"
$ws.Range("E15").Value = "  +0.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.23"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("E17").Value = "  -0.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008839"
$ws.Range("E18").Value = "  -1.55%  "

$ws.Range("E19").Value = "  -0.82%  "

$ws.Range("D20").Value = "28.155.09"
$ws.Range("E20").Value = "  +2.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.85"
$ws.Range("E21").Value = "  -0.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.163"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").Value = "2.177.38"
$ws.Range("E23").Value = "  +3.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  +1.87%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.51"
$ws.Range("E25").Value = "  +1.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.920"
$ws.Range("E26").Value = "  -1.78%  "

$ws.Range("E27").Value = "  -0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.104"
$ws.Range("E28").Value = "  +3.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.40"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.966"
$ws.Range("E30").Value = "  -0.77%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08913"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.366"
$ws.Range("E32").Value = "  +1.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.253"
$ws.Range("E33").Value = "  +3.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7717"
$ws.Range("E34").Value = "  +2.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.693"
$ws.Range("E35").Value = "  +1.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.647"
$ws.Range("E36").Value = "  -1.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02048"
$ws.Range("E37").Value = "  -1.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.103"
$ws.Range("E38").Value = "  -1.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05324"
$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5523"
$ws.Range("E40").Value = "  +2.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.000"
$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.067"
$ws.Range("E42").Value = "  -0.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1528"
$ws.Range("E43").Value = "  +0.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.470"
$ws.Range("E44").Value = "  +0.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.69"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4858"
$ws.Range("E46").Value = "  +0.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "107.55"
$ws.Range("E47").Value = "  +4.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.657"
$ws.Range("E49").Value = "  -0.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.57"
$ws.Range("E50").Value = "  +2.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06104"
$ws.Range("E51").Value = "  -0.18%  "
